# Add 2022-Q4 data:
#  - Duplicate the existing "2020-Q4" sheet so its original data is preserved
#    in a new tab (archived, unchanged).
#  - Rename/overwrite the original "2020-Q4" sheet with the fresh 2022-Q4
#    fund-holdings data (it keeps its original position/rId, right after the
#    "总计" sheet).
#  - Update the "总计" (totals) sheet with a new row for 2022-Q4, pushing the
#    existing 2020-Q4 totals row down.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)
$q4old = $wb.Worksheets.Item(2)

# --- 1. Archive the current "2020-Q4" sheet as-is, right after itself. ---
$q4old.Copy($null, $q4old)
$archive = $wb.Worksheets.Item(3)
$archive.Name = "2020-Q4-archive-tmp"

# --- 2. Turn the original sheet into the new "2022-Q4" sheet. ---
$q4old.Name = "2022-Q4"
$q4old.Cells.Clear()

# Headers (style copied from the "总计" sheet's header style so the new
# sheet's formatting matches the existing convention).
$total.Range("B1").Copy($q4old.Range("B1:H1"))
$q4old.Range("B1").Value = "基金代码"
$q4old.Range("C1").Value = "基金名称"
$q4old.Range("D1").Value = "基金规模"
$q4old.Range("E1").Value = "股票总仓位"
$q4old.Range("F1").Value = "仓位占比"
$q4old.Range("G1").Value = "持有市值(亿元)"
$q4old.Range("H1").Value = "仓位排名"

# Row index column (A) style, same convention as "总计" sheet's A column.
$total.Range("A2").Copy($q4old.Range("A2:A8"))

$data = @(
    @(0, "014831", "兴银中证1000指数增强A", "1.40", "82.60", "2.04", "0.0286", 1),
    @(1, "014832", "兴银中证1000指数增强C", "1.01", "82.60", "2.04", "0.0206", 1),
    @(2, "850007", "海通智选一年持有期股票B", "0.30", "82.33", "0.66", "0.0020", 9),
    @(3, "004680", "前海开源裕瑞混合A", "0.31", "20.77", "0.55", "0.0017", 10),
    @(4, "850788", "海通智选一年持有期股票A", "0.20", "82.33", "0.66", "0.0013", 9),
    @(5, "006190", "前海开源裕瑞混合C", "0.10", "20.77", "0.55", "0.0006", 10),
    @(6, "850799", "海通智选一年持有期股票C", "0.00", "82.33", "0.66", 0, 9)
)

# Columns B-G hold text (fund code / name / ratios) in the source data,
# including numeric-looking strings like "1.40" or "014831" that must stay
# text (not get coerced into numbers, which would e.g. drop the leading
# zero). Force text entry via NumberFormat "@" over the whole block, then
# strip the number-format override afterwards so no stray style index is
# left on the cells (matches the original sheet's un-styled data cells).
$q4old.Range("B2:G8").NumberFormat = "@"

$r = 2
foreach ($row in $data) {
    $q4old.Cells.Item($r, 1).Value = $row[0]
    $q4old.Cells.Item($r, 2).Value = $row[1]
    $q4old.Cells.Item($r, 3).Value = $row[2]
    $q4old.Cells.Item($r, 4).Value = $row[3]
    $q4old.Cells.Item($r, 5).Value = $row[4]
    $q4old.Cells.Item($r, 6).Value = $row[5]
    $q4old.Cells.Item($r, 7).Value = $row[6]
    $q4old.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

$q4old.Range("B2:G8").ClearFormats()

# Last row's "持有市值(亿元)" is a true 0 numeric value (not the text "0.0000").
$q4old.Cells.Item(8, 7).Value = 0

# --- 3. Rename the archived duplicate back to "2020-Q4". ---
$archive.Name = "2020-Q4"

# --- 4. Update the "总计" sheet: insert the 2022-Q4 row above 2020-Q4. ---
$total.Range("A2:D2").Copy($total.Range("A3:D3"))
$total.Range("A3").Value = 1

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 7
$total.Range("D2").Value = 0.05
